$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.994.41'
$ws.Range("E2").Value = '  +5.88%  '
$ws.Range("D3").Value = '3.677.29'
$ws.Range("E3").Value = '  +17.80%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '598.94'
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("D6").Value = '183.79'
$ws.Range("E6").Value = '  +5.32%  '
$ws.Range("D7").Value = '3.674.46'
$ws.Range("E7").Value = '  +17.86%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +3.92%  '
$ws.Range("E10").Value = '  +6.68%  '
$ws.Range("E11").Value = '  +3.37%  '
$ws.Range("D12").Value = '0.499'
$ws.Range("E12").Value = '  +4.80%  '
$ws.Range("D13").Value = '40.47'
$ws.Range("E13").Value = '  +11.81%  '
$ws.Range("E14").Value = '  +5.15%  '
$ws.Range("D15").Value = '4.291.59'
$ws.Range("E15").Value = '  +17.87%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.683.62'
$ws.Range("E16").Value = '  +18.31%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '71.024.75'
$ws.Range("E17").Value = '  +6.01%  '
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("E19").Value = '  +6.55%  '
$ws.Range("D20").Value = '17.02'
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").Value = '513.95'
$ws.Range("E21").Value = '  +5.40%  '
$ws.Range("D22").Value = '9.20'
$ws.Range("E22").Value = '  +17.02%  '
$ws.Range("D23").Value = '0.743'
$ws.Range("E23").Value = '  +6.89%  '
$ws.Range("D24").Value = '87.53'
$ws.Range("E24").Value = '  +4.52%  '
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").Value = '  +10.75%  '
$ws.Range("D26").Value = '13.50'
$ws.Range("E26").Value = '  +4.60%  '
$ws.Range("D27").Value = '11.04'
$ws.Range("E27").Value = '  +8.37%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E29").Value = '  +10.02%  '
$ws.Range("D30").Value = '8.18'
$ws.Range("E30").Value = '  +1.82%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '2.78'
$ws.Range("E31").Value = '  +6.60%  '
$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = '0.0000111'
$ws.Range("E32").Value = '  +17.56%  '
$ws.Range("D33").Value = '31.54'
$ws.Range("E33").Value = '  +11.73%  '
$ws.Range("E34").Value = '  +3.24%  '
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("D36").Value = '6.11'
$ws.Range("E36").Value = '  +8.38%  '
$ws.Range("E37").Value = '  +7.20%  '
$ws.Range("E38").Value = '  +11.15%  '
$ws.Range("E39").Value = '  +9.59%  '
$ws.Range("D40").Value = '51.20'
$ws.Range("E40").Value = '  +3.90%  '
$ws.Range("E41").Value = '  +3.53%  '
$ws.Range("D42").Value = '45.07'
$ws.Range("E42").Value = '  -6.28%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").Value = '8.85'
$ws.Range("E43").Value = '  +6.19%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.135.79'
$ws.Range("E44").Value = '  +11.70%  '
$ws.Range("D45").Value = '417.54'
$ws.Range("E45").Value = '  +12.05%  '
$ws.Range("E46").Value = '  +4.17%  '
$ws.Range("D47").Value = '0.0370'
$ws.Range("E47").Value = '  +6.05%  '
$ws.Range("D48").Value = '28.38'
$ws.Range("E48").Value = '  +15.14%  '
$ws.Range("D49").Value = '137.50'
$ws.Range("E49").Value = '  +2.09%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '2.48'
$ws.Range("E51").Value = '  +11.61%  '
